# Normalize casing of the "State"/"Year" header labels and the
# Massachusetts/Texas row labels, then leave the selection on A4 (the
# last edited cell), matching the authored commit's data cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "state"
$ws.Range("B2").Value = "year"
$ws.Range("A3").Value = "Massachusetts"
$ws.Range("A4").Value = "Texas"

$ws.Range("A4").Select() | Out-Null
